$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 16.432
$ws.Range("D12").Value = -7.328999999999999
$ws.Range("E14").Value = 16.994
$ws.Range("E26").Value = 16.756
$ws.Range("E31").Value = 17.12
$ws.Range("D32").Value = -8.052000000000001
$ws.Range("E35").Value = 16.63
$ws.Range("D36").Value = -7.783999999999999
$ws.Range("E37").Value = 16.737
$ws.Range("D38").Value = -7.662000000000001
$ws.Range("E45").Value = 16.96
$ws.Range("D46").Value = -8.263000000000002
$ws.Range("D54").Value = -8.494
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("E57").Value = 16.643
$ws.Range("D67").Value = -7.467999999999999
$ws.Range("D69").Value = -7.555
$ws.Range("D72").Value = -7.595000000000001
$ws.Range("D91").Value = -7.279000000000001
$ws.Range("D99").Value = -8.019
$ws.Range("E100").Value = 16.749
$ws.Range("E102").Value = 16.669
